# Rebuild the player-name header row and the data grid for the
# "Top Speed (m-s)" trainings/pregame sheet: several new player columns are
# inserted (interleaved alphabetically) among the existing ones, and a new
# "Alexis Rainey" column is added right after the Date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final header order for columns B..W (column A stays "Date").
$headers = @(
    "Alexis Rainey", "Balduzzi", "Burns", "Curley", "Doyle", "Espona",
    "Ferriolo", "Hackman", "Holzman", "Hughes", "Johnson", "McCann",
    "McFadden", "Medico", "Myers", "Pla", "Reilly", "Rodrigo", "Streib",
    "Tollaksen", "Wasyliw", "Yanovich"
)

# Final data values for rows 2..6, columns B..W (same order as $headers).
# $null marks a cell that has no recorded value for that player/date.
$rowsData = @(
    @(5.8045, 5.3423, 5.9378, $null, 5.7778, 6.3801, 5.7334, 6.2534, 5.9823, $null, 6.1689, 6.12,   6.1645, 6.0423, 6.3023, 5.8712, $null,  6.1467, $null,  5.9978, 6.6312, 6.2045),
    @(6.4045, 5.8289, 6.5578, 5.8712, 6.4623, 6.3023, 6.1334, 6.4489, 6.1112, 6.3956, 6.4245, 6.5045, 6.1778, 7.0756, 6.6934, 6.6534, 6.0667, 6.3623, $null,  5.7245, 6.8023, 5.8978),
    @(6.1756, 5.0223, 5.9156, 6.3578, 6.1312, 6.6445, 6.3023, 6.3556, 6.7956, 5.8978, 5.9245, 5.9267, 6.1112, 6.6156, 6.5112, 5.9067, 6.08,   5.6712, $null,  5.5689, 6.4578, $null),
    @(6.1734, 5.5467, $null, 5.3467, 6.14,   6.3623, 6.5823, 6.4667, 6.1289, 5.88,   5.8934, 6.2645, 5.8667, 6.4134, $null,  6.1134, 5.9378, 5.9912, 5.5445, 5.5734, 6.6312, $null),
    @(5.8756, 5.5,    $null, 5.0956, 5.7112, 6.1112, 6.9289, 6.5467, 5.9156, 4.9467, 6.1645, 5.4045, 5.8378, 5.8045, $null,  5.2267, 5.6178, 4.8023, 4.4045, $null,  6.6489, $null)
)

# Columns B..W are (partly) new; give the whole header row the same look
# (bold, bordered, centered/top aligned) as the existing "Date" header cell
# before we fill in the names, by copying its format across.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:W1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Write the header row (columns B..W == indices 2..23).
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value2 = $headers[$i]
}

# Write the data rows (sheet rows 2..6, columns B..W).
for ($r = 0; $r -lt $rowsData.Length; $r++) {
    $rowVals = $rowsData[$r]
    $sheetRow = $r + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($sheetRow, $c + 2).Value2 = $rowVals[$c]
    }
}
